# Remove the trailing " of own" from the title of the two
# "Sequence diagram" slides, turning:
#   "... : 1 / 2 of own" -> "... : 1 / 2"
#   "... : 2 / 2 of own" -> "... : 2 / 2"

$p = $ppt.ActivePresentation

$slideIndexes = @(12, 13)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(1)
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text
    $suffix = " of own"
    $pos = $fullText.IndexOf($suffix)
    if ($pos -ge 0) {
        $startChar = $pos + 1
        $len = $suffix.Length
        $chars = $tr.Characters($startChar, $len)
        $chars.Text = ""
    }
}
